$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value2 = 2907.3333  # H4
$ws.Cells.Item(4, 9).Value2 = 2111  # I4
$ws.Cells.Item(4, 11).Value2 = 2111  # K4
$ws.Cells.Item(4, 13).Value2 = -1997  # M4
$ws.Cells.Item(32, 8).Value2 = 43332  # H32
$ws.Cells.Item(32, 9).Value2 = 30001  # I32
$ws.Cells.Item(32, 10).Value2 = 49997.5  # J32
$ws.Cells.Item(32, 11).Value2 = 30001  # K32
$ws.Cells.Item(32, 12).Value2 = 49997.5  # L32
$ws.Cells.Item(32, 13).Value2 = -29675  # M32
$ws.Cells.Item(32, 14).Value2 = -50649.5  # N32
$ws.Cells.Item(40, 8).Value2 = 4017.647  # H40
$ws.Cells.Item(40, 10).Value2 = 4353.5713  # J40
$ws.Cells.Item(40, 12).Value2 = 4353.5713  # L40
$ws.Cells.Item(40, 14).Value2 = -4703.5713  # N40
$ws.Cells.Item(100, 8).Value2 = 5244.5  # H100
$ws.Cells.Item(100, 9).Value2 = 5810.6  # I100
$ws.Cells.Item(100, 11).Value2 = 5810.6  # K100
$ws.Cells.Item(100, 13).Value2 = -5269.6  # M100
$ws.Cells.Item(135, 8).Value2 = 26075.2  # H135
$ws.Cells.Item(135, 10).Value2 = 26075.2  # J135
$ws.Cells.Item(135, 12).Value2 = 234676.8  # L135
$ws.Cells.Item(135, 14).Value2 = -239746.8  # N135
$ws.Cells.Item(138, 8).Value2 = 2903.4358  # H138
$ws.Cells.Item(138, 10).Value2 = 3178.1765  # J138
$ws.Cells.Item(138, 12).Value2 = 9534.529500000001  # L138
$ws.Cells.Item(138, 14).Value2 = -19814.5295  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 14292553  # H32
$ws.Cells.Item(32, 9).Value2 = 14292553  # I32
$ws.Cells.Item(32, 11).Value2 = 14292553  # K32
$ws.Cells.Item(32, 13).Value2 = -14292266  # M32
$ws.Cells.Item(46, 8).Value2 = 4450  # H46
$ws.Cells.Item(46, 9).Value2 = 3999  # I46
$ws.Cells.Item(46, 11).Value2 = 3999  # K46
$ws.Cells.Item(46, 13).Value2 = -3680  # M46
$ws.Cells.Item(110, 8).Value2 = 1922  # H110
$ws.Cells.Item(110, 9).Value2 = 1402.5  # I110
$ws.Cells.Item(110, 11).Value2 = 1402.5  # K110
$ws.Cells.Item(110, 13).Value2 = 642.5  # M110
$ws.Cells.Item(122, 8).Value2 = 1624.5  # H122
$ws.Cells.Item(122, 9).Value2 = 1500  # I122
$ws.Cells.Item(122, 11).Value2 = 4500  # K122
$ws.Cells.Item(122, 13).Value2 = -2050  # M122

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value2 = 941.35486  # H94
$ws.Cells.Item(94, 9).Value2 = 968.34485  # I94
$ws.Cells.Item(94, 11).Value2 = 968.34485  # K94
$ws.Cells.Item(94, 13).Value2 = -517.34485  # M94
$ws.Cells.Item(105, 8).Value2 = 2253.8948  # H105
$ws.Cells.Item(105, 9).Value2 = 1733  # I105
$ws.Cells.Item(105, 10).Value2 = 2632.7273  # J105
$ws.Cells.Item(105, 11).Value2 = 1733  # K105
$ws.Cells.Item(105, 12).Value2 = 2632.7273  # L105
$ws.Cells.Item(105, 13).Value2 = 14  # M105
$ws.Cells.Item(105, 14).Value2 = -6126.7273  # N105
$ws.Cells.Item(134, 8).Value2 = 44826.832  # H134
$ws.Cells.Item(134, 9).Value2 = 1828.0555  # I134
$ws.Cells.Item(134, 11).Value2 = 5484.166499999999  # K134
$ws.Cells.Item(134, 13).Value2 = -2949.166499999999  # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value2 = 3209.9092  # H7
$ws.Cells.Item(7, 9).Value2 = 128  # I7
$ws.Cells.Item(7, 10).Value2 = 8603.25  # J7
$ws.Cells.Item(7, 11).Value2 = 128  # K7
$ws.Cells.Item(7, 12).Value2 = 8603.25  # L7
$ws.Cells.Item(7, 13).Value2 = -15  # M7
$ws.Cells.Item(7, 14).Value2 = -8829.25  # N7
$ws.Cells.Item(31, 8).Value2 = 780063.6  # H31
$ws.Cells.Item(31, 9).Value2 = 1914  # I31
$ws.Cells.Item(31, 11).Value2 = 1914  # K31
$ws.Cells.Item(31, 13).Value2 = -1619  # M31
$ws.Cells.Item(34, 8).Value2 = 780063.6  # H34
$ws.Cells.Item(34, 9).Value2 = 1914  # I34
$ws.Cells.Item(34, 11).Value2 = 1914  # K34
$ws.Cells.Item(34, 13).Value2 = -1712  # M34
$ws.Cells.Item(42, 8).Value2 = 4059  # H42
$ws.Cells.Item(42, 10).Value2 = 0  # J42
$ws.Cells.Item(42, 12).Value2 = 0  # L42
$ws.Cells.Item(42, 14).ClearContents()  # N42
$ws.Cells.Item(62, 8).Value2 = 4712.5  # H62
$ws.Cells.Item(62, 10).Value2 = 4712.5  # J62
$ws.Cells.Item(62, 12).Value2 = 4712.5  # L62
$ws.Cells.Item(62, 14).Value2 = -5960.5  # N62
$ws.Cells.Item(65, 8).Value2 = 4712.5  # H65
$ws.Cells.Item(65, 10).Value2 = 4712.5  # J65
$ws.Cells.Item(65, 12).Value2 = 23562.5  # L65
$ws.Cells.Item(65, 14).Value2 = -29802.5  # N65
$ws.Cells.Item(97, 8).Value2 = 82331.336  # H97
$ws.Cells.Item(97, 10).Value2 = 82331.336  # J97
$ws.Cells.Item(97, 12).Value2 = 82331.336  # L97
$ws.Cells.Item(97, 14).Value2 = -84313.336  # N97
$ws.Cells.Item(105, 8).Value2 = 1897.3  # H105
$ws.Cells.Item(105, 9).Value2 = 1581.6  # I105
$ws.Cells.Item(105, 11).Value2 = 1581.6  # K105
$ws.Cells.Item(105, 13).Value2 = 165.4000000000001  # M105
$ws.Cells.Item(122, 8).Value2 = 3118.8462  # H122
$ws.Cells.Item(122, 9).Value2 = 3045.4167  # I122
$ws.Cells.Item(122, 11).Value2 = 9136.250100000001  # K122
$ws.Cells.Item(122, 13).Value2 = -6686.250100000001  # M122

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value2 = 679181.4  # H12
$ws.Cells.Item(12, 9).Value2 = 425.2  # I12
$ws.Cells.Item(12, 10).Value2 = 1056268.1  # J12
$ws.Cells.Item(12, 11).Value2 = 1275.6  # K12
$ws.Cells.Item(12, 12).Value2 = 3168804.3  # L12
$ws.Cells.Item(12, 13).Value2 = -1102.6  # M12
$ws.Cells.Item(12, 14).Value2 = -3169150.3  # N12
$ws.Cells.Item(51, 8).Value2 = 9388.362999999999  # H51
$ws.Cells.Item(51, 9).Value2 = 5263  # I51
$ws.Cells.Item(51, 11).Value2 = 15789  # K51
$ws.Cells.Item(51, 13).Value2 = -15329  # M51
$ws.Cells.Item(56, 8).Value2 = 9999  # H56
$ws.Cells.Item(56, 9).Value2 = 9999  # I56
$ws.Cells.Item(56, 11).Value2 = 9999  # K56
$ws.Cells.Item(56, 13).Value2 = -9469  # M56
$ws.Cells.Item(132, 8).Value2 = 1964.2941  # H132
$ws.Cells.Item(132, 9).Value2 = 2203.7778  # I132
$ws.Cells.Item(132, 10).Value2 = 1694.875  # J132
$ws.Cells.Item(132, 11).Value2 = 19834.0002  # K132
$ws.Cells.Item(132, 12).Value2 = 15253.875  # L132
$ws.Cells.Item(132, 13).Value2 = -17304.0002  # M132
$ws.Cells.Item(132, 14).Value2 = -20313.875  # N132
$ws.Cells.Item(137, 8).Value2 = 4795.1875  # H137
$ws.Cells.Item(137, 10).Value2 = 3812.25  # J137
$ws.Cells.Item(137, 12).Value2 = 11436.75  # L137
$ws.Cells.Item(137, 14).Value2 = -21636.75  # N137

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value2 = 0  # H29
$ws.Cells.Item(29, 9).Value2 = 0  # I29
$ws.Cells.Item(29, 11).Value2 = 0  # K29
$ws.Cells.Item(29, 13).ClearContents()  # M29
$ws.Cells.Item(70, 8).Value2 = 5000  # H70
$ws.Cells.Item(70, 9).Value2 = 5000  # I70
$ws.Cells.Item(70, 11).Value2 = 5000  # K70
$ws.Cells.Item(70, 13).Value2 = -4730  # M70
$ws.Cells.Item(73, 8).Value2 = 5000  # H73
$ws.Cells.Item(73, 9).Value2 = 5000  # I73
$ws.Cells.Item(73, 11).Value2 = 5000  # K73
$ws.Cells.Item(73, 13).Value2 = -4064  # M73
$ws.Cells.Item(80, 8).Value2 = 3575  # H80
$ws.Cells.Item(80, 10).Value2 = 3750  # J80
$ws.Cells.Item(80, 12).Value2 = 3750  # L80
$ws.Cells.Item(80, 14).Value2 = -5746  # N80
$ws.Cells.Item(83, 8).Value2 = 3575  # H83
$ws.Cells.Item(83, 10).Value2 = 3750  # J83
$ws.Cells.Item(83, 12).Value2 = 18750  # L83
$ws.Cells.Item(83, 14).Value2 = -28734  # N83
$ws.Cells.Item(132, 8).Value2 = 125003560  # H132
$ws.Cells.Item(132, 9).Value2 = 166669920  # I132
$ws.Cells.Item(132, 10).Value2 = 4500  # J132
$ws.Cells.Item(132, 11).Value2 = 500009760  # K132
$ws.Cells.Item(132, 12).Value2 = 13500  # L132
$ws.Cells.Item(132, 13).Value2 = -500007230  # M132
$ws.Cells.Item(132, 14).Value2 = -18560  # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value2 = 49503  # H20
$ws.Cells.Item(20, 9).Value2 = 19000  # I20
$ws.Cells.Item(20, 10).Value2 = 80006  # J20
$ws.Cells.Item(20, 11).Value2 = 19000  # K20
$ws.Cells.Item(20, 12).Value2 = 80006  # L20
$ws.Cells.Item(20, 13).Value2 = -18774  # M20
$ws.Cells.Item(20, 14).Value2 = -80458  # N20
$ws.Cells.Item(22, 8).Value2 = 4188  # H22
$ws.Cells.Item(22, 9).Value2 = 4188  # I22
$ws.Cells.Item(22, 11).Value2 = 4188  # K22
$ws.Cells.Item(22, 13).Value2 = -3893  # M22
$ws.Cells.Item(26, 8).Value2 = 10010  # H26
$ws.Cells.Item(26, 9).Value2 = 0  # I26
$ws.Cells.Item(26, 11).Value2 = 0  # K26
$ws.Cells.Item(26, 13).ClearContents()  # M26
$ws.Cells.Item(27, 8).Value2 = 4188  # H27
$ws.Cells.Item(27, 9).Value2 = 4188  # I27
$ws.Cells.Item(27, 11).Value2 = 4188  # K27
$ws.Cells.Item(27, 13).Value2 = -4081  # M27
$ws.Cells.Item(46, 8).Value2 = 3436.1155  # H46
$ws.Cells.Item(46, 9).Value2 = 3202.1052  # I46
$ws.Cells.Item(46, 11).Value2 = 3202.1052  # K46
$ws.Cells.Item(46, 13).Value2 = -3014.1052  # M46
$ws.Cells.Item(55, 8).Value2 = 58823930  # H55
$ws.Cells.Item(55, 9).Value2 = 66667076  # I55
$ws.Cells.Item(55, 10).Value2 = 316.5  # J55
$ws.Cells.Item(55, 11).Value2 = 66667076  # K55
$ws.Cells.Item(55, 12).Value2 = 316.5  # L55
$ws.Cells.Item(55, 13).Value2 = -66666903  # M55
$ws.Cells.Item(55, 14).Value2 = -662.5  # N55
$ws.Cells.Item(56, 8).Value2 = 19816.334  # H56
$ws.Cells.Item(56, 9).Value2 = 9633.666999999999  # I56
$ws.Cells.Item(56, 11).Value2 = 9633.666999999999  # K56
$ws.Cells.Item(56, 13).Value2 = -8942.666999999999  # M56
$ws.Cells.Item(93, 8).Value2 = 83335680  # H93
$ws.Cells.Item(93, 9).Value2 = 111112910  # I93
$ws.Cells.Item(93, 11).Value2 = 111112910  # K93
$ws.Cells.Item(93, 13).Value2 = -111111662  # M93
$ws.Cells.Item(100, 8).Value2 = 0  # H100
$ws.Cells.Item(100, 9).Value2 = 0  # I100
$ws.Cells.Item(100, 11).Value2 = 0  # K100
$ws.Cells.Item(100, 13).ClearContents()  # M100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value2 = 1249.5  # H81
$ws.Cells.Item(81, 9).Value2 = 1266.3334  # I81
$ws.Cells.Item(81, 10).Value2 = 1199  # J81
$ws.Cells.Item(81, 11).Value2 = 2532.6668  # K81
$ws.Cells.Item(81, 12).Value2 = 2398  # L81
$ws.Cells.Item(81, 13).Value2 = -1471.6668  # M81
$ws.Cells.Item(81, 14).Value2 = -4520  # N81
$ws.Cells.Item(84, 8).Value2 = 1249.5  # H84
$ws.Cells.Item(84, 9).Value2 = 1266.3334  # I84
$ws.Cells.Item(84, 10).Value2 = 1199  # J84
$ws.Cells.Item(84, 11).Value2 = 12663.334  # K84
$ws.Cells.Item(84, 12).Value2 = 11990  # L84
$ws.Cells.Item(84, 13).Value2 = -7359.333999999999  # M84
$ws.Cells.Item(84, 14).Value2 = -22598  # N84
$ws.Cells.Item(100, 8).Value2 = 1305.8  # H100
$ws.Cells.Item(100, 9).Value2 = 1305.8  # I100
$ws.Cells.Item(100, 11).Value2 = 2611.6  # K100
$ws.Cells.Item(100, 13).Value2 = -2070.6  # M100
